$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.792.85'
$ws.Range('E2').Value = '  -2.40%  '

$ws.Range('D3').Value = '1.616.62'
$ws.Range('E3').Value = '  -2.24%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.48%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.15'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.35%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.40%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3944'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.67%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3835'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.87%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.002'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.45%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '49.47'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.46%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.349'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.52%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08452'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.09%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '23.63'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -5.57%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.039'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.44%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.584'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.33%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001280'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.01%  '

$ws.Range('D17').Value = '1.616.51'
$ws.Range('E17').Value = '  -2.71%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '93.70'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.67%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06930'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.63%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.91'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -6.23%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.814'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.62%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.54%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.40'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.90%  '

$ws.Range('D24').Value = '23.808.26'
$ws.Range('E24').Value = '  -2.40%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.484'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +5.98%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.820'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.54%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.21'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.31%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '156.81'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.32%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '140.01'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.25%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.296'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -8.37%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.748'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.71%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.492'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.13%  '

$ws.Range('D33').Value = '1.800.72'
$ws.Range('E33').Value = '  -2.55%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08091'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.41%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9812'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.27%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.616'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.36%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02872'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.30%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2663'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.94%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.09139'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.51%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '10.33'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.48%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '13.59'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.14%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.424'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.84%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.7506'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.68%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.02'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.04%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6913'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.41%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.463'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.94%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.068'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.09%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.001'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.39%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.08226'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.85%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '135.23'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.39%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.196'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -8.61%  '

